$wb = $excel.ActiveWorkbook

# --- Sheet1: rename Sheet1 -> TestCases and fill in the two test-case rows ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TestCases"

$ws1.Range("A1").Value = "TestCases"
$ws1.Range("B1").Value = "Runmode"

# --- Sheet2: create TestData right after TestCases by duplicating TestCases
#     (this avoids the extra default sheetFormatPr attribute a brand-new,
#     empty sheet created with Worksheets.Add() would otherwise pick up) ---
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "TestData"
$ws2.Cells.Clear() | Out-Null

# Finish filling in TestCases (column A top-to-bottom, then column B)
$ws1.Range("A2").Value = "AddCustomerTest"
$ws1.Range("A3").Value = "OpenAccountTest"
$ws1.Range("B2").Value = "Y"
$ws1.Range("B3").Value = "N"

$ws1.Columns.Item(1).ColumnWidth = 16.16666666666667
$ws1.Columns.Item(2).ColumnWidth = 8.666666666666666

# --- Fill in TestData: two stacked data-provider blocks, one per test case ---
$ws2.Range("A1").Value = "AddCustomerTest"

$ws2.Range("A2").Value = "Runmode"
$ws2.Range("B2").Value = "firstname"
$ws2.Range("C2").Value = "lastname"
$ws2.Range("D2").Value = "postcode"

$ws2.Range("A3").Value = "Y"
$ws2.Range("B3").Value = "manish"
$ws2.Range("C3").Value = "k"
$ws2.Range("D3").Value = "P6767"

$ws2.Range("A4").Value = "N"
$ws2.Range("C4").Value = "k"
$ws2.Range("D4").Value = "X7878"

$ws2.Range("A6").Value = "OpenAccountTest"

$ws2.Range("A7").Value = "Runmode"
$ws2.Range("B7").Value = "customer"
$ws2.Range("C7").Value = "currency"

$ws2.Range("A8").Value = "Y"
$ws2.Range("B8").Value = "manish k"

$ws2.Range("B4").Value = "jyoti"
$ws2.Range("B9").Value = "jyoti k"

$ws2.Range("C8").Value = "Rupee"
$ws2.Range("C9").Value = "Dollar"

$ws2.Range("A9").Value = "Y"

# --- Selections matching the target sheetViews ---
$ws1.Range("A3").Select() | Out-Null
$ws2.Range("A6").Select() | Out-Null

# Make TestData the active / tab-selected sheet (matches activeTab="1")
$ws2.Activate() | Out-Null

Write-Host "done"
